{"js": "// Locate the paragraph that ends with the \"Sammendrag\" instructional text\n// (\"... og hva denne rapporten handler om\") inside the summary table cell,\n// then append two new \"Comment\"-styled paragraphs after it: a blank spacer\n// paragraph and a paragraph holding the actual Train Dispatch System summary.\nconst searchResults = context.document.body.search(\n  \"og hva denne rapporten handler om\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find the target 'Sammendrag' paragraph text.\");\n}\n\nconst targetParagraph = searchResults.items[0].paragraphs.getFirst();\nconst insertionRange = targetParagraph.getRange(\"End\");\nawait context.sync();\n\n// Build the OOXML for the two new paragraphs exactly as authored:\n//  1) an empty paragraph using the \"Comment\" style (blank spacer line)\n//  2) a \"Comment\" style paragraph with the summary text, where the actual\n//     content runs override italics/color back to plain black text, and the\n//     word \"Dispatch\" is wrapped in proofErr spell-check markers.\nconst newParagraphsXml =\n  '<w:p><w:pPr><w:pStyle w:val=\"Comment\"/></w:pPr></w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"Comment\"/><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\">Train </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t>Dispatch</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\"> System er et program som lagrer og h\u00e5ndterer togavganger p\u00e5 en enkelt stasjon. Programmet holder orden p\u00e5 tognummer, avgangstid, forsinkelser, spor, linjenavn og destinasjon. </w:t></w:r>' +\n  \"</w:p>\";\n\nconst packageXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  \"<pkg:xmlData>\" +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  \"</Relationships>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newParagraphsXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionRange.insertOoxml(packageXml, \"After\");\nawait context.sync();\n", "ps1": "# Locate the paragraph that ends with the \"Sammendrag\" instructional text\n# (\"... og hva denne rapporten handler om\") inside the summary table cell,\n# then append two new \"Comment\"-styled paragraphs after it: a blank spacer\n# paragraph and a paragraph holding the actual Train Dispatch System summary.\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$findRange.Find.Text = \"og hva denne rapporten handler om\"\n$found = $findRange.Find.Execute()\nif (-not $found) {\n    throw \"Could not find the target 'Sammendrag' paragraph text.\"\n}\n\n$targetParagraph = $findRange.Paragraphs(1)\n$insertionPoint = $d.Range($targetParagraph.Range.End, $targetParagraph.Range.End)\n\n# Build the OOXML for the two new paragraphs exactly as authored:\n#  1) an empty paragraph using the \"Comment\" style (blank spacer line)\n#  2) a \"Comment\" style paragraph with the summary text, where the actual\n#     content runs override italics/color back to plain black text (both the\n#     literal RGB value and the \"text1\" theme reference), and the word\n#     \"Dispatch\" is wrapped in proofErr spell-check markers.\n$newParagraphsXml = (\n  '<w:p><w:pPr><w:pStyle w:val=\"Comment\"/></w:pPr></w:p>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"Comment\"/><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\">Train </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t>Dispatch</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:i w:val=\"0\"/><w:iCs/><w:color w:val=\"000000\" w:themeColor=\"text1\"/></w:rPr><w:t xml:space=\"preserve\"> System er et program som lagrer og h\u00e5ndterer togavganger p\u00e5 en enkelt stasjon. Programmet holder orden p\u00e5 tognummer, avgangstid, forsinkelser, spor, linjenavn og destinasjon. </w:t></w:r>' +\n  '</w:p>'\n)\n\n$packageXml = (\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + $newParagraphsXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n)\n\n$insertionPoint.InsertXML($packageXml)\n"}
